# Applies the "added ifo gdp component analysis preprocessing" update:
# recalculated averages in the staircase matrix plus one newly appended
# diagonal of values (column K for row 11 down through column B for row 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - revised average
$ws.Range("K10").Value = 1.388747888886706

# Row 11 - revised average + newly appended value
$ws.Range("J11").Value = 0.444773652920949
$ws.Range("K11").Value = 0.2348700177716323

# Row 12 - revised average + newly appended value
$ws.Range("I12").Value = 0.4487415504340581
$ws.Range("J12").Value = 0.2388379152847414

# Row 13 - revised average + newly appended value
$ws.Range("H13").Value = 0.5843816406042994
$ws.Range("I13").Value = 0.3744780054549828

# Row 14 - revised average + newly appended value
$ws.Range("G14").Value = 0.3435754587486348
$ws.Range("H14").Value = 0.1336718235993181

# Row 15 - revised average + newly appended value
$ws.Range("F15").Value = 0.2982442434965384
$ws.Range("G15").Value = 0.08834060834722172

# Row 16 - revised average + newly appended value
$ws.Range("E16").Value = 0.2313828215604846
$ws.Range("F16").Value = 0.02147918641116785

# Row 17 - revised average + newly appended value
$ws.Range("D17").Value = 0.201796619203768
$ws.Range("E17").Value = -0.00810701594554874

# Row 18 - revised average + newly appended value
$ws.Range("C18").Value = 0.1836459624741271
$ws.Range("D18").Value = -0.02625767267518964

# Row 19 - revised average + newly appended value
$ws.Range("B19").Value = 0.1656141382254278
$ws.Range("C19").Value = -0.04428949692388896

# Row 20 - newly appended value
$ws.Range("B20").Value = -0.09587373626955231
